# Auto-generated edit script applying the Mandragora_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific
# leve rows across the ALC, ARM, BSM, CUL, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 19827.445
$ws.Range("I21").Value = 19689.4
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 19689.4
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = -19221.4
$ws.Range("N21").Value = -20936

$ws.Range("H23").Value = 19827.445
$ws.Range("I23").Value = 19689.4
$ws.Range("J23").Value = 20000
$ws.Range("K23").Value = 19689.4
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = -19455.4
$ws.Range("N23").Value = -20468

$ws.Range("H62").Value = 16949.521
$ws.Range("I62").Value = 4083.5454
$ws.Range("J62").Value = 300001
$ws.Range("K62").Value = 4083.5454
$ws.Range("L62").Value = 300001
$ws.Range("M62").Value = -3459.5454
$ws.Range("N62").Value = -301249

$ws.Range("H65").Value = 16949.521
$ws.Range("I65").Value = 4083.5454
$ws.Range("J65").Value = 300001
$ws.Range("K65").Value = 20417.727
$ws.Range("L65").Value = 1500005
$ws.Range("M65").Value = -17297.727
$ws.Range("N65").Value = -1506245

$ws.Range("H137").Value = 3046.1765
$ws.Range("I137").Value = 3976.111
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 11928.333
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -9378.332999999999
$ws.Range("N137").Value = -11100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2757.3438
$ws.Range("I45").Value = 2078.4827
$ws.Range("J45").Value = 9319.666999999999
$ws.Range("K45").Value = 2078.4827
$ws.Range("L45").Value = 9319.666999999999
$ws.Range("M45").Value = -1701.4827
$ws.Range("N45").Value = -10073.667

$ws.Range("H63").Value = 14645.5
$ws.Range("I63").Value = 4156.6665
$ws.Range("K63").Value = 4156.6665
$ws.Range("M63").Value = -3470.6665

$ws.Range("H66").Value = 14645.5
$ws.Range("I66").Value = 4156.6665
$ws.Range("K66").Value = 20783.3325
$ws.Range("M66").Value = -17351.3325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 17400
$ws.Range("J76").Value = 17400
$ws.Range("L76").Value = 17400
$ws.Range("N76").Value = -18030

$ws.Range("H79").Value = 17400
$ws.Range("J79").Value = 17400
$ws.Range("L79").Value = 17400
$ws.Range("N79").Value = -19584

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1404.4615
$ws.Range("I5").Value = 391
$ws.Range("J5").Value = 3026
$ws.Range("K5").Value = 1173
$ws.Range("L5").Value = 9078
$ws.Range("M5").Value = -1061
$ws.Range("N5").Value = -9302

$ws.Range("H31").Value = 2900
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2900
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 8700
$ws.Range("N31").Value = -9276
$ws.Range("M31").ClearContents()

$ws.Range("H44").Value = 449.36365
$ws.Range("I44").Value = 360.75
$ws.Range("K44").Value = 1082.25
$ws.Range("M44").Value = -684.25

$ws.Range("H49").Value = 4167.1113
$ws.Range("J49").Value = 4167.1113
$ws.Range("L49").Value = 12501.3339
$ws.Range("N49").Value = -12813.3339

$ws.Range("H54").Value = 3221
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 3221
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 9663
$ws.Range("N54").Value = -10781
$ws.Range("M54").ClearContents()

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H74").Value = 2650
$ws.Range("I74").Value = 1250
$ws.Range("K74").Value = 3750
$ws.Range("M74").Value = -2689

$ws.Range("H77").Value = 2650
$ws.Range("I77").Value = 1250
$ws.Range("K77").Value = 11250
$ws.Range("M77").Value = -5946

$ws.Range("H94").Value = 2731.1667
$ws.Range("I94").Value = 750
$ws.Range("J94").Value = 3721.75
$ws.Range("K94").Value = 2250
$ws.Range("L94").Value = 11165.25
$ws.Range("M94").Value = -1574
$ws.Range("N94").Value = -12517.25

$ws.Range("H98").Value = 1076.8
$ws.Range("J98").Value = 1516.6666
$ws.Range("L98").Value = 4549.9998
$ws.Range("N98").Value = -7545.9998

$ws.Range("H100").Value = 3271.4285
$ws.Range("J100").Value = 3271.4285
$ws.Range("L100").Value = 9814.2855
$ws.Range("N100").Value = -11436.2855

$ws.Range("H101").Value = 8500
$ws.Range("J101").Value = 8500
$ws.Range("L101").Value = 25500
$ws.Range("N101").Value = -30368

$ws.Range("H107").Value = 413.72726
$ws.Range("J107").Value = 424.8
$ws.Range("L107").Value = 1274.4
$ws.Range("N107").Value = -5114.4

$ws.Range("H114").Value = 835.25
$ws.Range("I114").Value = 219.28572
$ws.Range("J114").Value = 1451.2142
$ws.Range("K114").Value = 657.85716
$ws.Range("L114").Value = 4353.642599999999
$ws.Range("M114").Value = 2596.14284
$ws.Range("N114").Value = -10861.6426

$ws.Range("H122").Value = 3407.8164
$ws.Range("I122").Value = 492.8
$ws.Range("J122").Value = 3564.5376
$ws.Range("K122").Value = 4435.2
$ws.Range("L122").Value = 32080.8384
$ws.Range("M122").Value = -1985.2
$ws.Range("N122").Value = -36980.8384

$ws.Range("H135").Value = 1404.4615
$ws.Range("I135").Value = 391
$ws.Range("J135").Value = 3026
$ws.Range("K135").Value = 3519
$ws.Range("L135").Value = 27234
$ws.Range("M135").Value = -984
$ws.Range("N135").Value = -32304

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 800.75
$ws.Range("I22").Value = 374.2857
$ws.Range("J22").Value = 1030.3846
$ws.Range("K22").Value = 374.2857
$ws.Range("L22").Value = 1030.3846
$ws.Range("M22").Value = -79.28570000000002
$ws.Range("N22").Value = -1620.3846

$ws.Range("H27").Value = 800.75
$ws.Range("I27").Value = 374.2857
$ws.Range("J27").Value = 1030.3846
$ws.Range("K27").Value = 374.2857
$ws.Range("L27").Value = 1030.3846
$ws.Range("M27").Value = -267.2857
$ws.Range("N27").Value = -1244.3846

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5566.241
$ws.Range("I126").Value = 7091.2383
$ws.Range("J126").Value = 1563.125
$ws.Range("K126").Value = 21273.7149
$ws.Range("L126").Value = 4689.375
$ws.Range("M126").Value = -18803.7149
$ws.Range("N126").Value = -9629.375
